$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "aragon" dimension (column G) is re-curated to use the shared SDMX
# reference-area dimension instead of a bespoke IAEST dimension, and its
# values now come from a "URI-Comunidad" typed column instead of a
# skos:Concept mapped through a dedicated mapping workbook.

# Row 2: dimension URI for column G
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 4: value class/type for column G
$ws.Range("G4").Value = "URI-Comunidad"

# Row 5: the per-dimension mapping workbook no longer applies to column G,
# so that cell is removed entirely (not just blanked).
$ws.Range("G5").Clear()
